$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ B = 0.4609249367666753; C = 0.7037408350246651; D = 0.6939948720570934; E = 0.833063546229874;  F = 0.7182889450858672; G = 15 }
    3  = @{ B = 0.3216291049761831; C = 0.6088832333399757; D = 0.584052044952158;  E = 0.7642329781893464; F = 0.7194279080246428; G = 14 }
    4  = @{ B = 0.2154567245497249; C = 0.5333824792787472; D = 0.4853342657017765; E = 0.6966593613106599; F = 0.6895569019380638; G = 13 }
    5  = @{ B = 0.4452680593528981; C = 0.5566372887181056; D = 0.4832252965807901; E = 0.6951440833243063; F = 0.5575547117719994; G = 12 }
    6  = @{ B = 0.3156011815745178; C = 0.4418649605082216; D = 0.2461552685847581; E = 0.4961403718553431; F = 0.4015050174666302; G = 11 }
    7  = @{ B = 0.3772995006961687; C = 0.4404536501459158; D = 0.2519683207200895; E = 0.5019644616106697; F = 0.348988072853359;  G = 10 }
    8  = @{ B = 0.3543917895297701; C = 0.4159301033722839; D = 0.2171525390636058; E = 0.4659962865341373; F = 0.3209421651944155; G = 9  }
    9  = @{ B = 0.380213756763938;  C = 0.4480587750649256; D = 0.2483796526124562; E = 0.4983770185436486; F = 0.3444534416474213; G = 8  }
    10 = @{ B = 0.3348409171772554; C = 0.4048669237504626; D = 0.2023859592503501; E = 0.4498732702110118; F = 0.3245182676007908; G = 7  }
    11 = @{ B = 0.3934609991409326; C = 0.465587436599884;  D = 0.2489390515396543; E = 0.4989379235332331; F = 0.3360847994682415; G = 6  }
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$row").Value = $rowVals[$col]
    }
}
